# Swap the species-observation data between row 4 and row 5,
# for columns A, B, E, F, G, H, Q, R (the other columns are identical
# between the two rows already, so they are left untouched).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A", "B", "E", "F", "G", "H", "Q", "R")

foreach ($col in $cols) {
    $addr4 = "${col}4"
    $addr5 = "${col}5"

    $val4 = $ws.Range($addr4).Value2
    $val5 = $ws.Range($addr5).Value2

    $ws.Range($addr4).Value = $val5
    $ws.Range($addr5).Value = $val4
}
